$d = $word.ActiveDocument

function Set-ParaText($para, [string]$text) {
    $r = $d.Range($para.Range.Start, $para.Range.End)
    $r.Text = $text
}

# ---------------------------------------------------------------------------
# Paragraph 2: "In the income statement..." -> "The two steps involved..."
#   + jc=both
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
Set-ParaText $p2 "The two steps involved in the physical counting of the inventory are as follows:"
$p2.Format.Alignment = 3

# ---------------------------------------------------------------------------
# Paragraph 3: "It includes ... buyer." -> list item "Planning" (bold)
#   pStyle=ListParagraph, numPr ilvl0/numId1, jc=both
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Style = "List Paragraph"
$full3 = $d.Range($p3.Range.Start, $p3.Range.End)
$full3.ListFormat.ApplyNumberDefault()
$p3.Format.Alignment = 3
$p3.Range.Font.Name = "Lucida Sans Typewriter"

$full3b = $d.Range($p3.Range.Start, $p3.Range.End)
$xml3 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/><w:b/><w:bCs/></w:rPr><w:t>Planning</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full3b.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Paragraph 4: "In the accounting perspective..." -> "Usually before..."
#   + ind left=720, jc=both
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
Set-ParaText $p4 "Usually before the inventory count begins, the firm must decide on a period during which the firm’s staff will involve in counting the inventory. "
$p4.Format.Alignment = 3
$p4.Format.LeftIndent = 36

# ---------------------------------------------------------------------------
# Paragraph 5: "COGS = Beginning Inventory..." (was center+bold) ->
#   "This is often done during the period..." (ind left=720, jc=both, no bold)
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
Set-ParaText $p5 "This is often done during the period when the firm experiences the lowest sales (often a trough in their business cycle!) and so, this is usually determined as the period when the counting occurs."
$p5.Format.Alignment = 3
$p5.Format.LeftIndent = 36
$full5 = $d.Range($p5.Range.Start, $p5.Range.End)
$full5.Font.Bold = 0
$full5.Font.Name = "Lucida Sans Typewriter"
$p5.Range.Font.Name = "Lucida Sans Typewriter"

# ---------------------------------------------------------------------------
# Paragraph 6: "The above equation is obvious..." -> list item
#   "Counting & Verification" (bold)
#   pStyle=ListParagraph, numPr ilvl0/numId1, jc=both
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$p6.Style = "List Paragraph"
$full6 = $d.Range($p6.Range.Start, $p6.Range.End)
$full6.ListFormat.ApplyNumberDefault()
$p6.Format.Alignment = 3
$p6.Range.Font.Name = "Lucida Sans Typewriter"

$full6b = $d.Range($p6.Range.Start, $p6.Range.End)
$xml6 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/><w:b/><w:bCs/></w:rPr><w:t>Counting &amp; Verification</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full6b.InsertXML($xml6)

# ---------------------------------------------------------------------------
# Paragraph 7: "This equation is used often..." -> "The second step involves..."
#   + ind left=720, jc=both (already had jc=both)
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
Set-ParaText $p7 "The second step involves counting the inventory manually by the staff. The inventory is calculated by counting and valuing the items and adding them all (if there are branches of the firm). "
$p7.Format.Alignment = 3
$p7.Format.LeftIndent = 36

# ---------------------------------------------------------------------------
# New paragraph 8: "After the counting is done..."
#   ind left=720, jc=both
# ---------------------------------------------------------------------------
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
Set-ParaText $p8 "After the counting is done, an external auditor can verify the count by bringing in experts and taking a small sample to verify the authenticity of the count."
$p8.Format.Alignment = 3
$p8.Format.LeftIndent = 36

# ---------------------------------------------------------------------------
# New paragraph 9: empty
#   ind left=720, jc=both
# ---------------------------------------------------------------------------
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs(9)
$p9.Format.Alignment = 3
$p9.Format.LeftIndent = 36

Write-Output "done; paragraphs=$($d.Paragraphs.Count)"
